$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.973.50"
$ws.Range("E2").Value = "  +1.40%  "
$ws.Range("D3").Value = "2.307.91"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "541.86"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").Value = "129.48"
$ws.Range("E6").Value = "  -2.23%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "0.573"
$ws.Range("E8").Value = "  -2.43%  "
$ws.Range("D9").Value = "2.306.10"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").Value = "0.101"
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("D11").Value = "5.56"
$ws.Range("E11").Value = "  +2.35%  "
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").Value = "0.332"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("E14").Value = "  -2.16%  "
$ws.Range("D15").Value = "59.930.25"
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("D16").Value = "2.719.14"
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "0.0000132"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("D18").Value = "2.340.21"
$ws.Range("E18").Value = "  +1.62%  "
$ws.Range("D19").Value = "10.49"
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("E20").Value = "  -2.09%  "
$ws.Range("D21").Value = "311.98"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").Value = "6.57"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +1.62%  "
$ws.Range("E25").Value = "  -1.65%  "
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("E27").Value = "  -2.59%  "
$ws.Range("D28").Value = "1.35"
$ws.Range("E28").Value = "  +4.06%  "
$ws.Range("D29").Value = "171.44"
$ws.Range("E29").Value = "  +0.89%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").Value = "0.0₃0726"
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("E34").Value = "  +3.10%  "
$ws.Range("D35").Value = "0.379"
$ws.Range("E35").Value = "  -1.16%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("D40").Value = "316.59"
$ws.Range("E40").Value = "  +2.41%  "
$ws.Range("D41").Value = "38.10"
$ws.Range("E41").Value = "  -0.50%  "
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("D43").Value = "136.03"
$ws.Range("E43").Value = "  -3.72%  "
$ws.Range("D44").Value = "3.42"
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("D45").Value = "0.0937"
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("D46").Value = "0.568"
$ws.Range("E46").Value = "  +2.06%  "
$ws.Range("D47").Value = "18.80"
$ws.Range("E47").Value = "  +2.13%  "
$ws.Range("D48").Value = "0.0490"
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("E49").Value = "  +22.04%  "
$ws.Range("E50").Value = "  +0.73%  "
$ws.Range("E51").Value = "  +0.14%  "
